$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.23420996532603056,
    -0.0099999991878298999,
    -0.0089999991817304448,
    0.061992493638598489,
    -0.0059999991965984378,
    -0.0059999991725518953,
    -0.019999999034221005,
    -0.019999999028780024,
    0.024337355267431526,
    -0.061013195588749625,
    -0.0044999991660361616,
    0.063773218056798253,
    -0.0059999991384822593,
    -0.011999999076333978,
    -0.00599999913189464,
    -0.0059999991292984944,
    -0.0059999991257440044,
    -0.0089999990959492848,
    -0.0089999991972544713,
    -0.0089999991897151688,
    -0.0089999991886777764,
    -0.0089999991880418406,
    -0.0089999991745219887,
    -0.041999998839770569,
    -0.041999998833591512,
    -0.0059999991705517175,
    -0.0059999991676527031,
    -0.0059999991563799426,
    -0.0119999990906976,
    -0.019999999009992386,
    0.047810648378444398,
    -0.057279068981669568,
    -0.0059999991399246611
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
